$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(3, 1).Value = 67914826
$ws.Cells.Item(3, 2).Value = 108194
$ws.Cells.Item(3, 5).Value = 219711
$ws.Cells.Item(3, 6).Value = "Sårläka"
$ws.Cells.Item(3, 7).Value = "Sanicula europaea"
$ws.Cells.Item(3, 8).Value = "L."
$ws.Cells.Item(3, 9).Value = ""
$ws.Cells.Item(3, 10).Value = ""
$ws.Cells.Item(3, 17).Value = 717052.9258912012
$ws.Cells.Item(3, 18).Value = 6646811.022708044
$ws.Cells.Item(3, 29).Value = ""
$ws.Cells.Item(4, 1).Value = 67914804
$ws.Cells.Item(4, 2).Value = 88896
$ws.Cells.Item(4, 4).Value = "VU"
$ws.Cells.Item(4, 5).Value = 720
$ws.Cells.Item(4, 6).Value = "Violgubbe"
$ws.Cells.Item(4, 7).Value = "Gomphus clavatus"
$ws.Cells.Item(4, 8).Value = "(Pers.) Gray"
$ws.Cells.Item(4, 9).Value = "7"
$ws.Cells.Item(4, 10).Value = "fruktkroppar"
$ws.Cells.Item(4, 17).Value = 717054.0436666304
$ws.Cells.Item(4, 18).Value = 6646791.97993397
$ws.Cells.Item(5, 1).Value = 67913258
$ws.Cells.Item(5, 2).Value = 90319
$ws.Cells.Item(5, 4).Value = "LC"
$ws.Cells.Item(5, 5).Value = 4769
$ws.Cells.Item(5, 6).Value = "Svavelriska"
$ws.Cells.Item(5, 7).Value = "Lactarius scrobiculatus"
$ws.Cells.Item(5, 8).Value = "(Scop.:Fr.) Fr."
$ws.Cells.Item(5, 9).Value = ""
$ws.Cells.Item(5, 10).Value = ""
$ws.Cells.Item(5, 17).Value = 716944.959965172
$ws.Cells.Item(5, 18).Value = 6646874.079515913
$ws.Cells.Item(6, 1).Value = 67913270
$ws.Cells.Item(6, 2).Value = 90661
$ws.Cells.Item(6, 4).Value = "VU"
$ws.Cells.Item(6, 5).Value = 2058
$ws.Cells.Item(6, 6).Value = "Koppartaggsvamp"
$ws.Cells.Item(6, 7).Value = "Hydnellum lundellii"
$ws.Cells.Item(6, 8).Value = "(Maas Geest. & Nannf.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Cells.Item(6, 9).Value = "2"
$ws.Cells.Item(6, 10).Value = "fruktkroppar"
$ws.Cells.Item(6, 17).Value = 716880.2087797265
$ws.Cells.Item(6, 18).Value = 6646878.829376921
$ws.Cells.Item(6, 29).Value = "I anslutning till vildsvinsbökad mark."
$ws.Cells.Item(7, 1).Value = 67913209
$ws.Cells.Item(7, 2).Value = 56411
$ws.Cells.Item(7, 4).Value = "NT"
$ws.Cells.Item(7, 5).Value = 100049
$ws.Cells.Item(7, 6).Value = "Spillkråka"
$ws.Cells.Item(7, 7).Value = "Dryocopus martius"
$ws.Cells.Item(7, 8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(7, 9).Value = "1"
$ws.Cells.Item(7, 13).Value = "födosökande"
$ws.Cells.Item(7, 17).Value = 717094.2343623195
$ws.Cells.Item(7, 18).Value = 6646948.209421237
$ws.Cells.Item(7, 19).Value = 50
$ws.Cells.Item(8, 1).Value = 67913282
$ws.Cells.Item(8, 2).Value = 89412
$ws.Cells.Item(8, 4).Value = "NT"
$ws.Cells.Item(8, 5).Value = 5442
$ws.Cells.Item(8, 6).Value = "Tallticka"
$ws.Cells.Item(8, 7).Value = "Porodaedalea pini"
$ws.Cells.Item(8, 8).Value = "(Brot.) Murrill"
$ws.Cells.Item(8, 13).Value = ""
$ws.Cells.Item(8, 17).Value = 716818.1044006473
$ws.Cells.Item(8, 18).Value = 6647053.18770974
$ws.Cells.Item(9, 1).Value = 67913278
$ws.Cells.Item(9, 2).Value = 4717
$ws.Cells.Item(9, 5).Value = 102306
$ws.Cells.Item(9, 6).Value = "Granbarkgnagare"
$ws.Cells.Item(9, 7).Value = "Microbregma emarginatum"
$ws.Cells.Item(9, 8).Value = "(Duftschmid, 1825)"
$ws.Cells.Item(9, 13).Value = "äldre gnagspår"
$ws.Cells.Item(9, 17).Value = 716755.7881303673
$ws.Cells.Item(9, 18).Value = 6646896.172787146
$ws.Cells.Item(10, 1).Value = 67913139
$ws.Cells.Item(10, 2).Value = 98520
$ws.Cells.Item(10, 5).Value = 222498
$ws.Cells.Item(10, 6).Value = "Blåsippa"
$ws.Cells.Item(10, 7).Value = "Hepatica nobilis"
$ws.Cells.Item(10, 8).Value = "Schreb."
$ws.Cells.Item(10, 17).Value = 717187.9496040216
$ws.Cells.Item(10, 18).Value = 6646939.130624352
$ws.Cells.Item(10, 29).Value = "Allmänt förekommande inom området."
$ws.Cells.Item(11, 1).Value = 67913851
$ws.Cells.Item(11, 2).Value = 90074
$ws.Cells.Item(11, 4).Value = "LC"
$ws.Cells.Item(11, 5).Value = 3298
$ws.Cells.Item(11, 6).Value = "Trådticka"
$ws.Cells.Item(11, 7).Value = "Climacocystis borealis"
$ws.Cells.Item(11, 8).Value = "(Fr.) Kotl. & Pouzar"
$ws.Cells.Item(11, 9).Value = ""
$ws.Cells.Item(11, 10).Value = ""
$ws.Cells.Item(11, 17).Value = 716689.1208222503
$ws.Cells.Item(11, 18).Value = 6646907.853172779
$ws.Cells.Item(11, 29).Value = "På höstubbe av gran."
$ws.Cells.Item(12, 1).Value = 67913822
$ws.Cells.Item(12, 2).Value = 90674
$ws.Cells.Item(12, 4).Value = "LC"
$ws.Cells.Item(12, 5).Value = 5964
$ws.Cells.Item(12, 6).Value = "Fjällig taggsvamp s.str."
$ws.Cells.Item(12, 7).Value = "Sarcodon imbricatus s.str."
$ws.Cells.Item(12, 8).Value = "(L.:Fr.) P.Karst."
$ws.Cells.Item(12, 9).Value = ""
$ws.Cells.Item(12, 13).Value = ""
$ws.Cells.Item(12, 17).Value = 716671.0292412415
$ws.Cells.Item(12, 18).Value = 6646932.940137157
$ws.Cells.Item(12, 19).Value = 5
$ws.Cells.Item(12, 29).Value = "Allmänt förekommande inom området."
$ws.Cells.Item(13, 1).Value = 67913289
$ws.Cells.Item(13, 2).Value = 73631
$ws.Cells.Item(13, 5).Value = 6426
$ws.Cells.Item(13, 6).Value = "Kattfotslav"
$ws.Cells.Item(13, 7).Value = "Felipes leucopellaeus"
$ws.Cells.Item(13, 8).Value = "(Ach.) Frisch & G.Thor"
$ws.Cells.Item(13, 17).Value = 716779.0899203173
$ws.Cells.Item(13, 18).Value = 6647014.194535435
$ws.Cells.Item(13, 29).Value = "På gammal gran."
$ws.Cells.Item(14, 1).Value = 67914797
$ws.Cells.Item(14, 2).Value = 89410
$ws.Cells.Item(14, 4).Value = "NT"
$ws.Cells.Item(14, 5).Value = 5432
$ws.Cells.Item(14, 6).Value = "Granticka"
$ws.Cells.Item(14, 7).Value = "Porodaedalea chrysoloma"
$ws.Cells.Item(14, 8).Value = "(Fr.) Fiasson & Niemelä"
$ws.Cells.Item(14, 17).Value = 717063.1636913288
$ws.Cells.Item(14, 18).Value = 6646773.909642578
$ws.Cells.Item(14, 29).Value = "På grenar på gammal levande gran."
$ws.Cells.Item(15, 1).Value = 67914512
$ws.Cells.Item(15, 2).Value = 90005
$ws.Cells.Item(15, 4).Value = "LC"
$ws.Cells.Item(15, 5).Value = 1339
$ws.Cells.Item(15, 6).Value = "Brandticka"
$ws.Cells.Item(15, 7).Value = "Pycnoporellus fulgens"
$ws.Cells.Item(15, 8).Value = "(Fr.) Donk"
$ws.Cells.Item(15, 17).Value = 716658.238045899
$ws.Cells.Item(15, 18).Value = 6646558.06836888
$ws.Cells.Item(15, 29).Value = "På grov högstubbe av gran."
$ws.Cells.Item(16, 1).Value = 67913188
$ws.Cells.Item(16, 2).Value = 90665
$ws.Cells.Item(16, 5).Value = 4366
$ws.Cells.Item(16, 6).Value = "Skarp dropptaggsvamp"
$ws.Cells.Item(16, 7).Value = "Hydnellum peckii"
$ws.Cells.Item(16, 8).Value = "Banker"
$ws.Cells.Item(16, 13).Value = ""
$ws.Cells.Item(16, 17).Value = 717179.1824814338
$ws.Cells.Item(16, 18).Value = 6646951.186559737
$ws.Cells.Item(17, 1).Value = 67913183
$ws.Cells.Item(17, 17).Value = 717196.8702048012
$ws.Cells.Item(17, 18).Value = 6646958.762340512
$ws.Cells.Item(18, 1).Value = 67913274
$ws.Cells.Item(18, 2).Value = 73631
$ws.Cells.Item(18, 5).Value = 6426
$ws.Cells.Item(18, 6).Value = "Kattfotslav"
$ws.Cells.Item(18, 7).Value = "Felipes leucopellaeus"
$ws.Cells.Item(18, 8).Value = "(Ach.) Frisch & G.Thor"
$ws.Cells.Item(18, 17).Value = 716880.97296897
$ws.Cells.Item(18, 18).Value = 6646865.800371076
$ws.Cells.Item(18, 29).Value = "På gammal gran."
$ws.Cells.Item(19, 1).Value = 67913066
$ws.Cells.Item(19, 2).Value = 101680
$ws.Cells.Item(19, 5).Value = 222412
$ws.Cells.Item(19, 6).Value = "Tibast"
$ws.Cells.Item(19, 7).Value = "Daphne mezereum"
$ws.Cells.Item(19, 8).Value = "L."
$ws.Cells.Item(19, 29).Value = ""
$ws.Cells.Item(20, 1).Value = 67913254
$ws.Cells.Item(20, 2).Value = 101680
$ws.Cells.Item(20, 5).Value = 222412
$ws.Cells.Item(20, 6).Value = "Tibast"
$ws.Cells.Item(20, 7).Value = "Daphne mezereum"
$ws.Cells.Item(20, 8).Value = "L."
$ws.Cells.Item(20, 17).Value = 716951.0444857149
$ws.Cells.Item(20, 18).Value = 6646881.979102565
$ws.Cells.Item(20, 29).Value = ""
$ws.Cells.Item(21, 1).Value = 67914880
$ws.Cells.Item(21, 2).Value = 90074
$ws.Cells.Item(21, 5).Value = 3298
$ws.Cells.Item(21, 6).Value = "Trådticka"
$ws.Cells.Item(21, 7).Value = "Climacocystis borealis"
$ws.Cells.Item(21, 8).Value = "(Fr.) Kotl. & Pouzar"
$ws.Cells.Item(21, 17).Value = 717163.0512021189
$ws.Cells.Item(21, 18).Value = 6646899.955352104
$ws.Cells.Item(21, 29).Value = "På torrgran."
$ws.Cells.Item(22, 1).Value = 67913211
$ws.Cells.Item(22, 17).Value = 717094.2343623195
$ws.Cells.Item(22, 18).Value = 6646948.209421237
$ws.Cells.Item(23, 1).Value = 67914856
$ws.Cells.Item(23, 2).Value = 90665
$ws.Cells.Item(23, 5).Value = 4366
$ws.Cells.Item(23, 6).Value = "Skarp dropptaggsvamp"
$ws.Cells.Item(23, 7).Value = "Hydnellum peckii"
$ws.Cells.Item(23, 8).Value = "Banker"
$ws.Cells.Item(23, 17).Value = 717129.134103466
$ws.Cells.Item(23, 18).Value = 6646799.908343226
$ws.Cells.Item(23, 29).Value = ""
$ws.Cells.Item(24, 1).Value = 67913818
$ws.Cells.Item(24, 2).Value = 90665
$ws.Cells.Item(24, 5).Value = 4366
$ws.Cells.Item(24, 6).Value = "Skarp dropptaggsvamp"
$ws.Cells.Item(24, 7).Value = "Hydnellum peckii"
$ws.Cells.Item(24, 8).Value = "Banker"
$ws.Cells.Item(24, 17).Value = 716671.0292412415
$ws.Cells.Item(24, 18).Value = 6646932.940137157
$ws.Cells.Item(24, 29).Value = ""
$ws.Cells.Item(25, 1).Value = 67914777
$ws.Cells.Item(25, 17).Value = 716852.0415224443
$ws.Cells.Item(25, 18).Value = 6646672.013723582
$ws.Cells.Item(26, 1).Value = 67913202
$ws.Cells.Item(26, 2).Value = 89410
$ws.Cells.Item(26, 4).Value = "NT"
$ws.Cells.Item(26, 5).Value = 5432
$ws.Cells.Item(26, 6).Value = "Granticka"
$ws.Cells.Item(26, 7).Value = "Porodaedalea chrysoloma"
$ws.Cells.Item(26, 8).Value = "(Fr.) Fiasson & Niemelä"
$ws.Cells.Item(26, 17).Value = 717132.8554776106
$ws.Cells.Item(26, 18).Value = 6646950.979990459
$ws.Cells.Item(26, 29).Value = "På grenar på torrgran."
$ws.Cells.Item(27, 1).Value = 67914522
$ws.Cells.Item(27, 2).Value = 98520
$ws.Cells.Item(27, 5).Value = 222498
$ws.Cells.Item(27, 6).Value = "Blåsippa"
$ws.Cells.Item(27, 7).Value = "Hepatica nobilis"
$ws.Cells.Item(27, 8).Value = "Schreb."
$ws.Cells.Item(27, 29).Value = "Allmänt förekommande inom området."
$ws.Cells.Item(28, 1).Value = 67913247
$ws.Cells.Item(28, 2).Value = 89410
$ws.Cells.Item(28, 4).Value = "NT"
$ws.Cells.Item(28, 5).Value = 5432
$ws.Cells.Item(28, 6).Value = "Granticka"
$ws.Cells.Item(28, 7).Value = "Porodaedalea chrysoloma"
$ws.Cells.Item(28, 8).Value = "(Fr.) Fiasson & Niemelä"
$ws.Cells.Item(28, 17).Value = 716979.9164224314
$ws.Cells.Item(28, 18).Value = 6646982.22850688
$ws.Cells.Item(28, 29).Value = "På grenar på torrgran nära kärr."
$ws.Cells.Item(29, 1).Value = 67913143
$ws.Cells.Item(29, 2).Value = 90697
$ws.Cells.Item(29, 4).Value = "NT"
$ws.Cells.Item(29, 5).Value = 5449
$ws.Cells.Item(29, 6).Value = "Svart taggsvamp"
$ws.Cells.Item(29, 7).Value = "Phellodon niger"
$ws.Cells.Item(29, 8).Value = "(Fr.:Fr.) P.Karst."
$ws.Cells.Item(29, 17).Value = 717185.8272155134
$ws.Cells.Item(29, 18).Value = 6646958.113689412
$ws.Cells.Item(30, 1).Value = 73101105
$ws.Cells.Item(30, 2).Value = 85077
$ws.Cells.Item(30, 5).Value = 3762
$ws.Cells.Item(30, 6).Value = "Olivspindling"
$ws.Cells.Item(30, 7).Value = "Cortinarius venetus"
$ws.Cells.Item(30, 8).Value = "(Fr.:Fr.) Fr."
$ws.Cells.Item(30, 9).Value = "2"
$ws.Cells.Item(30, 10).Value = "fruktkroppar"
$ws.Cells.Item(30, 16).Value = "Storkärret, Södra Rörvik, Upl"
$ws.Cells.Item(30, 17).Value = 716788.3853775742
$ws.Cells.Item(30, 18).Value = 6646795.504157254
$ws.Cells.Item(30, 19).Value = 100
$ws.Cells.Item(30, 25).Value = "2018-09-09"
$ws.Cells.Item(30, 27).Value = "2018-09-09"
$ws.Cells.Item(30, 49).Value = "Jacob Rudhe"
$ws.Cells.Item(30, 50).Value = "Jacob Rudhe, Rasmus Elleby, Jon Jörpeland"
$ws.Cells.Item(31, 1).Value = 73101093
$ws.Cells.Item(31, 2).Value = 96355
$ws.Cells.Item(31, 5).Value = 219862
$ws.Cells.Item(31, 6).Value = "Nästrot"
$ws.Cells.Item(31, 7).Value = "Neottia nidus-avis"
$ws.Cells.Item(31, 8).Value = "(L.) Rich."
$ws.Cells.Item(31, 9).Value = "1"
$ws.Cells.Item(31, 10).Value = "stjälkar/strån/skott"
$ws.Cells.Item(31, 11).Value = "överblommad"
$ws.Cells.Item(31, 16).Value = "Storkärret, Södra Rörvik, Upl"
$ws.Cells.Item(31, 17).Value = 716677.9603077322
$ws.Cells.Item(31, 18).Value = 6646548.160545002
$ws.Cells.Item(31, 19).Value = 50
$ws.Cells.Item(31, 25).Value = "2018-09-09"
$ws.Cells.Item(31, 27).Value = "2018-09-09"
$ws.Cells.Item(31, 29).Value = ""
$ws.Cells.Item(31, 49).Value = "Jacob Rudhe"
$ws.Cells.Item(31, 50).Value = "Jacob Rudhe, Rasmus Elleby, Jon Jörpeland"
$ws.Cells.Item(32, 1).Value = 73101110
$ws.Cells.Item(32, 2).Value = 96312
$ws.Cells.Item(32, 5).Value = 219798
$ws.Cells.Item(32, 6).Value = "Skogsknipprot"
$ws.Cells.Item(32, 7).Value = "Epipactis helleborine"
$ws.Cells.Item(32, 8).Value = "(L.) Crantz"
$ws.Cells.Item(32, 11).Value = "överblommad"
$ws.Cells.Item(32, 16).Value = "Storkärret, Södra Rörvik, Upl"
$ws.Cells.Item(32, 17).Value = 717050.8085189572
$ws.Cells.Item(32, 18).Value = 6646718.372707062
$ws.Cells.Item(32, 19).Value = 100
$ws.Cells.Item(32, 25).Value = "2018-09-09"
$ws.Cells.Item(32, 27).Value = "2018-09-09"
$ws.Cells.Item(32, 29).Value = ""
$ws.Cells.Item(32, 49).Value = "Jacob Rudhe"
$ws.Cells.Item(32, 50).Value = "Jacob Rudhe, Rasmus Elleby, Jon Jörpeland"
$ws.Cells.Item(33, 1).Value = 73101095
$ws.Cells.Item(33, 2).Value = 108194
$ws.Cells.Item(33, 5).Value = 219711
$ws.Cells.Item(33, 6).Value = "Sårläka"
$ws.Cells.Item(33, 7).Value = "Sanicula europaea"
$ws.Cells.Item(33, 11).Value = "överblommad"
$ws.Cells.Item(33, 16).Value = "Storkärret, Södra Rörvik, Upl"
$ws.Cells.Item(33, 17).Value = 716677.9603077322
$ws.Cells.Item(33, 18).Value = 6646548.160545002
$ws.Cells.Item(33, 19).Value = 50
$ws.Cells.Item(33, 25).Value = "2018-09-09"
$ws.Cells.Item(33, 27).Value = "2018-09-09"
$ws.Cells.Item(33, 49).Value = "Jacob Rudhe"
$ws.Cells.Item(33, 50).Value = "Jacob Rudhe, Rasmus Elleby, Jon Jörpeland"
$ws.Cells.Item(34, 1).Value = 111667309
$ws.Cells.Item(34, 2).Value = 5113
$ws.Cells.Item(34, 5).Value = 100526
$ws.Cells.Item(34, 6).Value = "Bronshjon"
$ws.Cells.Item(34, 7).Value = "Callidium coriaceum"
$ws.Cells.Item(34, 8).Value = "Paykull, 1800"
$ws.Cells.Item(34, 9).Value = "1"
$ws.Cells.Item(34, 13).Value = "äldre gnagspår"
$ws.Cells.Item(34, 16).Value = "Södra Rörvik, Upl"
$ws.Cells.Item(34, 17).Value = 717045
$ws.Cells.Item(34, 18).Value = 6646924
$ws.Cells.Item(34, 19).Value = 100
$ws.Cells.Item(34, 25).Value = "2023-08-24"
$ws.Cells.Item(34, 26).Value = ""
$ws.Cells.Item(34, 27).Value = "2023-08-24"
$ws.Cells.Item(34, 28).Value = ""
$ws.Cells.Item(34, 29).Value = "Larvgångar och ingångshål i veden på många granar som varit döda i några år."
$ws.Cells.Item(34, 49).Value = "Åke Lindelöw"
$ws.Cells.Item(34, 50).Value = "Åke Lindelöw"
$ws.Cells.Item(35, 1).Value = 67914010
$ws.Cells.Item(35, 2).Value = 93235
$ws.Cells.Item(35, 5).Value = 210
$ws.Cells.Item(35, 6).Value = "Grön sköldmossa"
$ws.Cells.Item(35, 7).Value = "Buxbaumia viridis"
$ws.Cells.Item(35, 8).Value = "(Moug. ex Lam. & DC.) Brid. ex Moug. & Nestl."
$ws.Cells.Item(35, 9).Value = "1"
$ws.Cells.Item(35, 10).Value = "kapslar"
$ws.Cells.Item(35, 17).Value = 716545.1521927882
$ws.Cells.Item(35, 18).Value = 6646588.154634802
$ws.Cells.Item(35, 29).Value = "På gammal rötskadad granlåga."
$ws.Cells.Item(36, 1).Value = 67914177
$ws.Cells.Item(36, 2).Value = 108194
$ws.Cells.Item(36, 5).Value = 219711
$ws.Cells.Item(36, 6).Value = "Sårläka"
$ws.Cells.Item(36, 7).Value = "Sanicula europaea"
$ws.Cells.Item(36, 17).Value = 716551.2037129711
$ws.Cells.Item(36, 18).Value = 6646502.014087326
$ws.Cells.Item(37, 1).Value = 67914236
$ws.Cells.Item(37, 2).Value = 5113
$ws.Cells.Item(37, 5).Value = 100526
$ws.Cells.Item(37, 6).Value = "Bronshjon"
$ws.Cells.Item(37, 7).Value = "Callidium coriaceum"
$ws.Cells.Item(37, 8).Value = "Paykull, 1800"
$ws.Cells.Item(37, 13).Value = "äldre gnagspår"
$ws.Cells.Item(37, 17).Value = 716606.9416252106
$ws.Cells.Item(37, 18).Value = 6646462.030048631
$ws.Cells.Item(37, 29).Value = ""
$ws.Cells.Item(38, 1).Value = 67914140
$ws.Cells.Item(38, 2).Value = 85077
$ws.Cells.Item(38, 5).Value = 3762
$ws.Cells.Item(38, 6).Value = "Olivspindling"
$ws.Cells.Item(38, 7).Value = "Cortinarius venetus"
$ws.Cells.Item(38, 8).Value = "(Fr.:Fr.) Fr."
$ws.Cells.Item(38, 17).Value = 716558.9255464342
$ws.Cells.Item(38, 18).Value = 6646542.193822215
$ws.Cells.Item(38, 29).Value = ""
$ws.Cells.Item(39, 1).Value = 67914480
$ws.Cells.Item(39, 17).Value = 716678.2056377706
$ws.Cells.Item(39, 18).Value = 6646397.810981502
$ws.Cells.Item(40, 1).Value = 67914073
$ws.Cells.Item(40, 2).Value = 99398
$ws.Cells.Item(40, 5).Value = 221235
$ws.Cells.Item(40, 6).Value = "Vårärt"
$ws.Cells.Item(40, 7).Value = "Lathyrus vernus"
$ws.Cells.Item(40, 8).Value = "(L.) Bernh."
$ws.Cells.Item(40, 17).Value = 716558.9255464342
$ws.Cells.Item(40, 18).Value = 6646542.193822215
$ws.Cells.Item(41, 1).Value = 67914379
$ws.Cells.Item(41, 17).Value = 716677.1987222674
$ws.Cells.Item(41, 18).Value = 6646389.202686813
$ws.Cells.Item(42, 1).Value = 67914145
$ws.Cells.Item(42, 2).Value = 90674
$ws.Cells.Item(42, 4).Value = "LC"
$ws.Cells.Item(42, 5).Value = 5964
$ws.Cells.Item(42, 6).Value = "Fjällig taggsvamp s.str."
$ws.Cells.Item(42, 7).Value = "Sarcodon imbricatus s.str."
$ws.Cells.Item(42, 8).Value = "(L.:Fr.) P.Karst."
$ws.Cells.Item(42, 17).Value = 716558.9255464342
$ws.Cells.Item(42, 18).Value = 6646542.193822215
$ws.Cells.Item(42, 29).Value = "Allmänt förekommande inom området."
$ws.Cells.Item(43, 1).Value = 67914377
$ws.Cells.Item(43, 2).Value = 99398
$ws.Cells.Item(43, 4).Value = "LC"
$ws.Cells.Item(43, 5).Value = 221235
$ws.Cells.Item(43, 6).Value = "Vårärt"
$ws.Cells.Item(43, 7).Value = "Lathyrus vernus"
$ws.Cells.Item(43, 8).Value = "(L.) Bernh."
$ws.Cells.Item(43, 17).Value = 716677.1987222674
$ws.Cells.Item(43, 18).Value = 6646389.202686813
$ws.Cells.Item(43, 29).Value = ""
$ws.Cells.Item(44, 1).Value = 67914156
$ws.Cells.Item(44, 2).Value = 90665
$ws.Cells.Item(44, 5).Value = 4366
$ws.Cells.Item(44, 6).Value = "Skarp dropptaggsvamp"
$ws.Cells.Item(44, 7).Value = "Hydnellum peckii"
$ws.Cells.Item(44, 8).Value = "Banker"
$ws.Cells.Item(44, 17).Value = 716540.1688511835
$ws.Cells.Item(44, 18).Value = 6646527.014893715
$ws.Cells.Item(44, 29).Value = ""
$ws.Cells.Item(45, 1).Value = 67914325
$ws.Cells.Item(45, 2).Value = 73631
$ws.Cells.Item(45, 5).Value = 6426
$ws.Cells.Item(45, 6).Value = "Kattfotslav"
$ws.Cells.Item(45, 7).Value = "Felipes leucopellaeus"
$ws.Cells.Item(45, 8).Value = "(Ach.) Frisch & G.Thor"
$ws.Cells.Item(45, 17).Value = 716625.9801276408
$ws.Cells.Item(45, 18).Value = 6646438.000615025
$ws.Cells.Item(45, 29).Value = "På gammal gran."
$ws.Cells.Item(46, 1).Value = 67914183
$ws.Cells.Item(46, 2).Value = 90319
$ws.Cells.Item(46, 4).Value = "LC"
$ws.Cells.Item(46, 5).Value = 4769
$ws.Cells.Item(46, 6).Value = "Svavelriska"
$ws.Cells.Item(46, 7).Value = "Lactarius scrobiculatus"
$ws.Cells.Item(46, 8).Value = "(Scop.:Fr.) Fr."
$ws.Cells.Item(46, 17).Value = 716584.1544428288
$ws.Cells.Item(46, 18).Value = 6646489.862934983
$ws.Cells.Item(46, 29).Value = ""
$ws.Cells.Item(47, 1).Value = 67914364
$ws.Cells.Item(47, 2).Value = 90665
$ws.Cells.Item(47, 4).Value = "LC"
$ws.Cells.Item(47, 5).Value = 4366
$ws.Cells.Item(47, 6).Value = "Skarp dropptaggsvamp"
$ws.Cells.Item(47, 7).Value = "Hydnellum peckii"
$ws.Cells.Item(47, 8).Value = "Banker"
$ws.Cells.Item(47, 17).Value = 716675.1848726226
$ws.Cells.Item(47, 18).Value = 6646371.986039684
$ws.Cells.Item(48, 1).Value = 67914198
$ws.Cells.Item(48, 2).Value = 99398
$ws.Cells.Item(48, 5).Value = 221235
$ws.Cells.Item(48, 6).Value = "Vårärt"
$ws.Cells.Item(48, 7).Value = "Lathyrus vernus"
$ws.Cells.Item(48, 8).Value = "(L.) Bernh."
$ws.Cells.Item(48, 16).Value = "Norra Rörvik, S om, Upl"
$ws.Cells.Item(48, 17).Value = 716578.2384365106
$ws.Cells.Item(48, 18).Value = 6646461.857673807
$ws.Cells.Item(48, 19).Value = 5
$ws.Cells.Item(48, 25).Value = "2017-09-28"
$ws.Cells.Item(48, 27).Value = "2017-09-28"
$ws.Cells.Item(48, 49).Value = "Bo Törnquist"
$ws.Cells.Item(48, 50).Value = "Bo Törnquist, Kjell  Andersson"
$ws.Cells.Item(49, 1).Value = 67914209
$ws.Cells.Item(49, 2).Value = 98520
$ws.Cells.Item(49, 5).Value = 222498
$ws.Cells.Item(49, 6).Value = "Blåsippa"
$ws.Cells.Item(49, 7).Value = "Hepatica nobilis"
$ws.Cells.Item(49, 8).Value = "Schreb."
$ws.Cells.Item(49, 9).Value = ""
$ws.Cells.Item(49, 10).Value = ""
$ws.Cells.Item(49, 16).Value = "Norra Rörvik, S om, Upl"
$ws.Cells.Item(49, 17).Value = 716578.2384365106
$ws.Cells.Item(49, 18).Value = 6646461.857673807
$ws.Cells.Item(49, 19).Value = 5
$ws.Cells.Item(49, 25).Value = "2017-09-28"
$ws.Cells.Item(49, 27).Value = "2017-09-28"
$ws.Cells.Item(49, 29).Value = "Allmänt förekommande inom området."
$ws.Cells.Item(49, 49).Value = "Bo Törnquist"
$ws.Cells.Item(49, 50).Value = "Bo Törnquist, Kjell  Andersson"
$ws.Cells.Item(50, 1).Value = 67914162
$ws.Cells.Item(50, 2).Value = 99398
$ws.Cells.Item(50, 5).Value = 221235
$ws.Cells.Item(50, 6).Value = "Vårärt"
$ws.Cells.Item(50, 7).Value = "Lathyrus vernus"
$ws.Cells.Item(50, 8).Value = "(L.) Bernh."
$ws.Cells.Item(50, 9).Value = ""
$ws.Cells.Item(50, 10).Value = ""
$ws.Cells.Item(50, 11).Value = ""
$ws.Cells.Item(50, 16).Value = "Norra Rörvik, S om, Upl"
$ws.Cells.Item(50, 17).Value = 716532.9289578204
$ws.Cells.Item(50, 18).Value = 6646513.013206748
$ws.Cells.Item(50, 19).Value = 5
$ws.Cells.Item(50, 25).Value = "2017-09-28"
$ws.Cells.Item(50, 27).Value = "2017-09-28"
$ws.Cells.Item(50, 49).Value = "Bo Törnquist"
$ws.Cells.Item(50, 50).Value = "Bo Törnquist, Kjell  Andersson"
$ws.Cells.Item(51, 1).Value = 67914178
$ws.Cells.Item(51, 2).Value = 101680
$ws.Cells.Item(51, 5).Value = 222412
$ws.Cells.Item(51, 6).Value = "Tibast"
$ws.Cells.Item(51, 7).Value = "Daphne mezereum"
$ws.Cells.Item(51, 8).Value = "L."
$ws.Cells.Item(51, 11).Value = ""
$ws.Cells.Item(51, 16).Value = "Norra Rörvik, S om, Upl"
$ws.Cells.Item(51, 17).Value = 716551.2037129711
$ws.Cells.Item(51, 18).Value = 6646502.014087326
$ws.Cells.Item(51, 19).Value = 5
$ws.Cells.Item(51, 25).Value = "2017-09-28"
$ws.Cells.Item(51, 27).Value = "2017-09-28"
$ws.Cells.Item(51, 49).Value = "Bo Törnquist"
$ws.Cells.Item(51, 50).Value = "Bo Törnquist, Kjell  Andersson"
$ws.Cells.Item(52, 1).Value = 67914185
$ws.Cells.Item(52, 2).Value = 90674
$ws.Cells.Item(52, 5).Value = 5964
$ws.Cells.Item(52, 6).Value = "Fjällig taggsvamp s.str."
$ws.Cells.Item(52, 7).Value = "Sarcodon imbricatus s.str."
$ws.Cells.Item(52, 8).Value = "(L.:Fr.) P.Karst."
$ws.Cells.Item(52, 11).Value = ""
$ws.Cells.Item(52, 16).Value = "Norra Rörvik, S om, Upl"
$ws.Cells.Item(52, 17).Value = 716584.1544428288
$ws.Cells.Item(52, 18).Value = 6646489.862934983
$ws.Cells.Item(52, 19).Value = 5
$ws.Cells.Item(52, 25).Value = "2017-09-28"
$ws.Cells.Item(52, 27).Value = "2017-09-28"
$ws.Cells.Item(52, 29).Value = "Allmänt förekommande inom området."
$ws.Cells.Item(52, 49).Value = "Bo Törnquist"
$ws.Cells.Item(52, 50).Value = "Bo Törnquist, Kjell  Andersson"
$ws.Cells.Item(53, 1).Value = 73074510
$ws.Cells.Item(53, 2).Value = 90665
$ws.Cells.Item(53, 5).Value = 4366
$ws.Cells.Item(53, 6).Value = "Skarp dropptaggsvamp"
$ws.Cells.Item(53, 7).Value = "Hydnellum peckii"
$ws.Cells.Item(53, 8).Value = "Banker"
$ws.Cells.Item(53, 9).Value = ""
$ws.Cells.Item(53, 13).Value = ""
$ws.Cells.Item(53, 16).Value = "Södra Rörvik, SV om, Upl"
$ws.Cells.Item(53, 17).Value = 716692.486062364
$ws.Cells.Item(53, 18).Value = 6646360.427102973
$ws.Cells.Item(53, 19).Value = 10
$ws.Cells.Item(53, 25).Value = "2018-09-09"
$ws.Cells.Item(53, 26).Value = "00:00"
$ws.Cells.Item(53, 27).Value = "2018-09-09"
$ws.Cells.Item(53, 28).Value = "00:00"
$ws.Cells.Item(53, 29).Value = ""
$ws.Cells.Item(53, 49).Value = "Rasmus Elleby"
$ws.Cells.Item(53, 50).Value = "Rasmus Elleby, Jacob Rudhe, Jon Jörpeland"
